$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Recruiter Specialist"
$ws.Range("B2").Value = "For Hr"

$ws.Range("B9").Select()
